$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.510.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.633.72"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.114.33"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.296.33"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.634.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.09"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.05%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.34"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "548.03"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.23"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0301"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.82"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.92%  "
